# Fix credit_account in supplier management:
# The invoice in row 2 was associated with the wrong supplier / tax number.
# Update the supplier name ("المورد") and tax number ("الرقم الضريبي") for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: G = المورد (supplier), H = الرقم الضريبي (tax number)
$ws.Range("G2").Value = "مورد1"
$ws.Range("H2").Value = 3100000000003

# Leave the selection where the user last clicked before saving.
$ws.Range("H3").Select()
